$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 299, pushing existing rows 299:313 down to 300:314
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the new weekly data point.
# The "constant" columns (A,B,C,E,F,G,H,N,O,Q,R) mirror the rest of this block.
$ws.Range("A299").Value = 6
$ws.Range("B299").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C299").Value = "Metropolitana"
$ws.Range("D299").Value = 44509
$ws.Range("E299").Value = 13
$ws.Range("F299").Value = 100112039
$ws.Range("G299").Value = "Ciboulette"
$ws.Range("H299").Value = "Sin especificar"
$ws.Range("I299").Value = "Primera"
$ws.Range("J299").Value = 840
$ws.Range("K299").Value = 700
$ws.Range("L299").Value = 800
$ws.Range("M299").Value = 757
$ws.Range("N299").Value = "`$/docena de atados"
$ws.Range("O299").Value = "Región Metropolitana"
$ws.Range("P299").Value = 252
$ws.Range("Q299").Value = 3
$ws.Range("R299").Value = "Hortaliza"
